# Repull data, push all data, mean calculation
# Update the dSF (column F) values for several rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F8").Value = -1
$ws.Range("F13").Value = -1
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = -2
$ws.Range("F19").Value = -2
